$d = $word.ActiveDocument

# Locate the anchor paragraph: "This is a mash of opencv, posenet and pytorch. ..."
$anchor = $d.Content
$found = $anchor.Find.Execute("This is a mash of opencv, posenet and pytorch. I am going to give this a try.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "anchor paragraph not found"
}

# The paragraph object containing the found text.
$anchorPara = $anchor.Paragraphs(1)

# Insert a new empty paragraph right after the anchor paragraph; this becomes
# paragraph #1 of the new content ("Story short, ...").
$anchorPara.Range.InsertParagraphAfter()

# --- Paragraph 1: "Story short, I did not get it to work. ..." ---
$p1 = $anchorPara.Next()
$p1.Range.InsertAfter("Story short,")
$p1.Range.InsertAfter(" I")
$p1.Range.InsertAfter(" did not get it to work. I have tried it on windows, Linux wsl and an linux virtualmachine. I noticed the code is pretty outdated. It had not received an update in 2 years and when installing the recommended library it was not available anymore. ")

# New paragraph after p1 for paragraph 2.
$p1.Range.InsertParagraphAfter()

# --- Paragraph 2: "With this I realized ... NodeJS. I am using the npm serialport package(<link>). This is the only package ..." ---
$p2 = $p1.Next()
$p2.Range.InsertAfter("With this I ")
$p2.Range.InsertAfter("realized")
$p2.Range.InsertAfter(" I need to stop wasting time on something new and start with something I know, javascript. The first step for me is to get the serial communication working on ")
$p2.Range.InsertAfter("NodeJS")
$p2.Range.InsertAfter(".")
$p2.Range.InsertAfter(" I am using the npm serialport package(")

# Insert the hyperlink display text, then convert that range into a real hyperlink.
# (Range.End always points just past the paragraph mark, so subtract 1 to land
# on the actual text boundary before/after the inserted link text.)
$linkStart = $p2.Range.End - 1
$p2.Range.InsertAfter("https://duckduckgo.com/?t=ffab&q=nodejs+serial&ia=web")
$linkEnd = $p2.Range.End - 1
$hlRange = $d.Range($linkStart, $linkEnd)
$d.Hyperlinks.Add($hlRange, "https://duckduckgo.com/?t=ffab&q=nodejs+serial&ia=web", "", "", "https://duckduckgo.com/?t=ffab&q=nodejs+serial&ia=web")

$p2.Range.InsertAfter(" ). This is the only package nodejs as for javascript as the client side is still experimental and unavailable.")

# New paragraph after p2 for paragraph 3.
$p2.Range.InsertParagraphAfter()

# --- Paragraph 3: "The demo code can be found in demo_nodejs_communication." ---
$p3 = $p2.Next()
$p3.Range.InsertAfter("The demo code can be found in demo_nodejs_communication.")

# New paragraph after p3 for paragraph 4.
$p3.Range.InsertParagraphAfter()

# --- Paragraph 4: "Now this is working the next step is to get the data the robot needs" ---
$p4 = $p3.Next()
$p4.Range.InsertAfter("Now this is working the next step is to get the data the robot needs")

# New empty paragraph after p4 (the 5th added paragraph, kept blank).
$p4.Range.InsertParagraphAfter()

Write-Output "done"
